# Adonnai's Hours Log 3/20/15
# Add two new internship-log entries (rows 11 & 12) to Sheet1, matching
# the style/format of the existing log rows, and update the view state
# to reflect where the user ended up after typing the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11: Tue 3/11/15, 2:28 PM - 5:05 PM -------------------------------
$ws.Cells.Item(11, 1).Value = 0.60277777777777775   # Time In  (14:28)
$ws.Cells.Item(11, 2).Value = 0.71180555555555547   # Time Out (17:05)

# Copy the Date cell's number format down from the row above, then set
# the new date (3/11/2015) as a real value so it keeps the date style.
[void]$ws.Cells.Item(10, 3).Copy()
[void]$ws.Cells.Item(11, 3).PasteSpecial(-4122)      # xlPasteFormats
$ws.Cells.Item(11, 3).Value = 42074

# Duration column already carries the shared formula down through row 50.

# Copy the Notes cell's format from a later (already-styled) blank row so
# the new note picks up the same wrap/alignment style, then fill the text.
[void]$ws.Cells.Item(13, 5).Copy()
[void]$ws.Cells.Item(11, 5).PasteSpecial(-4122)      # xlPasteFormats
$ws.Cells.Item(11, 5).Value = "Modified the resizing function mechanics located in the html file. This gave the website more versatility between different screen-sizes. Hypothetically, phones, monitors, and tablets should be able to view the website with ease. "

# --- Row 12: Mon 3/16/15, 3:01 PM - 5:15 PM -------------------------------
$ws.Cells.Item(12, 1).Value = 0.62569444444444444   # Time In  (15:01)
$ws.Cells.Item(12, 2).Value = 0.71875               # Time Out (17:15)

[void]$ws.Cells.Item(10, 3).Copy()
[void]$ws.Cells.Item(12, 3).PasteSpecial(-4122)
$ws.Cells.Item(12, 3).Value = 42079

[void]$ws.Cells.Item(13, 5).Copy()
[void]$ws.Cells.Item(12, 5).PasteSpecial(-4122)
$ws.Cells.Item(12, 5).Value = "Created a function that highlights text in red when a mouse hovers over it. "

# --- View state: user scrolled down and left the cursor on the newest note
[void]$ws.Activate()
[void]$ws.Range("E14").Select()
$excel.ActiveWindow.Zoom = 83
